$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 9999
$ws.Range("J26").Value = 9999
$ws.Range("L26").Value = 9999
$ws.Range("N26").Value = -10687
$ws.Range("H28").Value = 1452.25
$ws.Range("J28").Value = 3526.3333
$ws.Range("L28").Value = 3526.3333
$ws.Range("N28").Value = -4496.3333
$ws.Range("H33").Value = 884.5714
$ws.Range("I33").Value = 884.5714
$ws.Range("K33").Value = 884.5714
$ws.Range("M33").Value = -655.5714
$ws.Range("H40").Value = 4559.1816
$ws.Range("J40").Value = 4706.375
$ws.Range("L40").Value = 4706.375
$ws.Range("N40").Value = -5056.375
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H98").Value = 2325.3333
$ws.Range("I98").Value = 2190.4
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 2190.4
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -692.4000000000001
$ws.Range("N98").Value = -5996
$ws.Range("H122").Value = 2325.3333
$ws.Range("I122").Value = 2190.4
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6571.200000000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4121.200000000001
$ws.Range("N122").Value = -13900
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3881.8572
$ws.Range("I32").Value = 2887
$ws.Range("K32").Value = 2887
$ws.Range("M32").Value = -2600
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H132").Value = 4328
$ws.Range("I132").Value = 4232
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 12696
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -10166
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2666.6667
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 2666.6667
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 2014
$ws.Range("I99").Value = 2014
$ws.Range("K99").Value = 2014
$ws.Range("M99").Value = -516
$ws.Range("H105").Value = 1616
$ws.Range("I105").Value = 1616
$ws.Range("K105").Value = 1616
$ws.Range("M105").Value = 131
$ws.Range("H134").Value = 537.5
$ws.Range("I134").Value = 537.5
$ws.Range("K134").Value = 1612.5
$ws.Range("M134").Value = 922.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 574.8333
$ws.Range("I16").Value = 574.8333
$ws.Range("K16").Value = 574.8333
$ws.Range("M16").Value = -287.8333
$ws.Range("H113").Value = 574.8333
$ws.Range("I113").Value = 574.8333
$ws.Range("K113").Value = 574.8333
$ws.Range("M113").Value = 1595.1667
$ws.Range("H134").Value = 4684.357
$ws.Range("I134").Value = 2533.7273
$ws.Range("J134").Value = 12570
$ws.Range("K134").Value = 7601.1819
$ws.Range("L134").Value = 37710
$ws.Range("M134").Value = -5066.1819
$ws.Range("N134").Value = -42780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 579.5714
$ws.Range("I2").Value = 636.7778
$ws.Range("J2").Value = 476.6
$ws.Range("K2").Value = 3820.6668
$ws.Range("L2").Value = 2859.6
$ws.Range("M2").Value = -3707.6668
$ws.Range("N2").Value = -3085.6
$ws.Range("H32").Value = 8300
$ws.Range("J32").Value = 8300
$ws.Range("L32").Value = 24900
$ws.Range("N32").Value = -25466
$ws.Range("H81").Value = 508.33334
$ws.Range("I81").Value = 210
$ws.Range("K81").Value = 630
$ws.Range("M81").Value = 493
$ws.Range("H84").Value = 508.33334
$ws.Range("I84").Value = 210
$ws.Range("K84").Value = 1890
$ws.Range("M84").Value = 3726
$ws.Range("H121").Value = 630.5
$ws.Range("J121").Value = 1022
$ws.Range("L121").Value = 3066
$ws.Range("N121").Value = -5686
$ws.Range("H131").Value = 1698
$ws.Range("I131").Value = 872.5
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 2617.5
$ws.Range("L131").Value = 15000
$ws.Range("M131").Value = 2422.5
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 470.26666
$ws.Range("I2").Value = 497.66666
$ws.Range("J2").Value = 360.66666
$ws.Range("K2").Value = 497.66666
$ws.Range("L2").Value = 360.66666
$ws.Range("M2").Value = -384.66666
$ws.Range("N2").Value = -586.66666
$ws.Range("H25").Value = 6479.8
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 7110.8887
$ws.Range("K25").Value = 800
$ws.Range("L25").Value = 7110.8887
$ws.Range("M25").Value = -271
$ws.Range("N25").Value = -8168.8887
$ws.Range("H31").Value = 2877
$ws.Range("I31").Value = 565.5
$ws.Range("K31").Value = 565.5
$ws.Range("M31").Value = -273.5
$ws.Range("H37").Value = 2877
$ws.Range("I37").Value = 565.5
$ws.Range("K37").Value = 565.5
$ws.Range("M37").Value = -288.5
$ws.Range("H132").Value = 2392.7778
$ws.Range("I132").Value = 2219.2856
$ws.Range("K132").Value = 6657.8568
$ws.Range("M132").Value = -4127.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2479.1667
$ws.Range("I100").Value = 2441.3635
$ws.Range("K100").Value = 2441.3635
$ws.Range("M100").Value = -1900.3635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 27054.084
$ws.Range("I32").Value = 12324.5
$ws.Range("K32").Value = 12324.5
$ws.Range("M32").Value = -12007.5
